# Automatische test-sync: 2025-06-30 20:10:50
#
# 1) "Logs" sheet: append a new log row (row 15) for testmail #15.
# 2) "Dashboard" sheet: swap the "Openingstijden / Locatie" and
#    "Bestelling / Levering" rows (row 4 <-> row 5), and append a new
#    "Overig" row (row 9).
# 3) Update the conditional-formatting ranges on "Logs" (D/G/H/I/J) so they
#    cover the newly added row.
# 4) Update the Dashboard bar chart's category/value series so it covers
#    the newly added Dashboard row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet - add row 15
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A15").Value = "Wat is jullie privacybeleid?"
$logs.Range("B15").Value = "mailmind.test@zohomail.eu"
$logs.Range("C15").Value = "Testmail #15: Wat is jullie privacybeleid?"
$logs.Range("D15").Value = "Overig"
$logs.Range("E15").Value = "Beste afzender,`nDank u voor uw vraag over ons privacybeleid. Wij hechten veel waarde aan de bescherming van persoonlijke gegevens van onze klanten en volgen strikte richtlijnen om deze te waarborgen. Ons privacybeleid is te vinden op onze website onder [link naar privacybeleid]. Hier vindt u gedetailleerde informatie over hoe wij omgaan met persoonlijke gegevens, welke gegevens wij verzamelen, hoe wij deze gebruiken en welke maatregelen wij treffen om ze te beschermen.`nMocht u nog verdere vragen hebben over ons privacybeleid, dan helpen wij graag verder.`nMet vriendelijke groet,`n[Naam] `nE-mailassistent"
$logs.Range("F15").Value = "2025-06-30 20:10:33"
$logs.Range("G15").Value = "Ja"
$logs.Range("H15").Value = "Nee"
$logs.Range("I15").Value = "Ja"
$logs.Range("J15").Value = "Nee"

# The multi-line Antwoord text auto-expands the row height; restore the
# default (non-custom) row height to match the other rows.
$logs.Rows.Item(15).AutoFit()

# ---------------------------------------------------------------------
# 2) Dashboard sheet - swap row 4 / row 5, add row 9
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Bestelling / Levering"
$dash.Range("B4").Value = 1
$dash.Range("A5").Value = "Openingstijden / Locatie"
$dash.Range("B5").Value = 1

$dash.Range("A9").Value = "Overig"
$dash.Range("B9").Value = 1

# ---------------------------------------------------------------------
# 3) Logs sheet - extend conditional formatting ranges to include row 15
# ---------------------------------------------------------------------
$ranges = "D", "G", "H", "I", "J"
foreach ($col in $ranges) {
    $oldRange = $logs.Range("$($col)2:$($col)14")
    $newRange = $logs.Range("$($col)2:$($col)15")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 4) Dashboard chart - extend category/value series to include row 9
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$9,'Dashboard'!`$B`$2:`$B`$9,1)"
